$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the two name values (B2, B3)
$ws.Range("B2").Value = "imron"
$ws.Range("B3").Value = "iksan"

# Apply a "Text" number format to column D data cells (D2:D3 and onward),
# which introduces a new style (numFmtId 49) without touching the header (D1)
$ws.Range("D2:D1048576").NumberFormat = "@"

# Hide columns E, F, G and J
$ws.Columns.Item(5).Hidden = $true
$ws.Columns.Item(6).Hidden = $true
$ws.Columns.Item(7).Hidden = $true
$ws.Columns.Item(10).Hidden = $true

# Update the selection to the whole column F
$ws.Range("F1:F1048576").Select()

Write-Host "done"
